# Updated ITA model - 2025-07-30 07:58
#
# Changes applied:
#  1. On the "System Settings" worksheet (sheet4.xml), insert a new row
#     above row 9 and populate it with a new attribute entry
#     (NCAP_AF / 0 / 3), pushing the rest of the table down by one row.
#  2. Make "System Settings" the active sheet/tab (was "fuels").
#  3. Select cell C11 on "System Settings" as the new active cell.

$wb = $excel.ActiveWorkbook

$sysSettings = $wb.Worksheets.Item("System Settings")

# Insert a new row above row 9 - everything below (rows 9-39) shifts down
# to rows 10-40, formatting is inherited from the row above as Excel does
# natively, and the sheet's used range/dimension grows to row 40.
$sysSettings.Rows("9:9").Insert()

# Populate the newly inserted row with the new attribute entry.
$sysSettings.Cells.Item(9, 3).Value = "NCAP_AF"
$sysSettings.Cells.Item(9, 4).Value = 0
$sysSettings.Cells.Item(9, 5).Value = 3

# Make "System Settings" the active sheet (this also clears tabSelected on
# the previously active "fuels" sheet automatically).
$sysSettings.Activate()

# Update the selection to match the new layout.
$sysSettings.Range("C11").Select()
